# Weekly update: insert a new price record (row 64) for Orégano at
# "Vega Central Mapocho de Santiago", pushing the existing historical
# rows (old 64..126) down by one to (65..127).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 64; this shifts rows 64:126 down to 65:127
# and Excel automatically grows the sheet dimension to A1:R127.
$ws.Rows("64:64").Insert()

# Populate the newly inserted row with the new week's data record.
$ws.Range("A64").Value = 9
$ws.Range("B64").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C64").Value = "Metropolitana"
$ws.Range("D64").Value = 45280
$ws.Range("E64").Value = 13
$ws.Range("F64").Value = 100112029
$ws.Range("G64").Value = "Orégano"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 16
$ws.Range("K64").Value = 21000
$ws.Range("L64").Value = 21000
$ws.Range("M64").Value = 21000
$ws.Range("N64").Value = "$/docena de atados"
$ws.Range("O64").Value = "Región Metropolitana"
$ws.Range("P64").Value = 7000
$ws.Range("Q64").Value = 3
$ws.Range("R64").Value = "Hortaliza"
